# Project upload feedback (Fixes #23)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the disease strain value in D2 to the full species name
$ws.Range("D2").Value = "Batrachochytridium dendrobatidus"

# Move the active selection to D3, matching the author's final cursor position
$ws.Range("D3").Select()
